$d = $word.ActiveDocument

# 1) "November 8, 2020" -> "November 13, 2020"
$d.Content.Find.Execute("November 8, 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "November 13, 2020", 2) | Out-Null

# 2) Company name in the body line: "at SpaceX." -> "at Palo Alto Networks."
$d.Content.Find.Execute("at SpaceX.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "at Palo Alto Networks.", 2) | Out-Null

# 3) Mission sentence rewrite
$d.Content.Find.Execute("align with SpaceX’s mission to craft the most advanced rockets and spacecraft.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "align with Palo Alto’s mission to deliver cutting-edge cybersecurity.", 2) | Out-Null
